$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 729, shifting existing rows 729:774 down to 730:775
$ws.Rows.Item(729).Insert()

# Populate the newly inserted row 729 with the new record's data
$ws.Range("A729").Value = 10
$ws.Range("B729").Value = "Vega Modelo de Temuco"
$ws.Range("C729").Value = "La Araucanía"
$ws.Range("D729").Value = 45041
$ws.Range("E729").Value = 9
$ws.Range("F729").Value = 100112043
$ws.Range("G729").Value = "Pepino ensalada"
$ws.Range("H729").Value = "Sin especificar"
$ws.Range("I729").Value = "Primera"
$ws.Range("J729").Value = 500
$ws.Range("K729").Value = 12000
$ws.Range("L729").Value = 14000
$ws.Range("M729").Value = 12800
$ws.Range("N729").Value = "$/caja 60 unidades"
$ws.Range("O729").Value = "Región de Arica y Parinacota"
$ws.Range("P729").Value = 213
$ws.Range("Q729").Value = 60
$ws.Range("R729").Value = "Hortaliza"
